# Atualização de bases das ligas, do dia: 13-06-2024 às 19:35
#
# The match rows in the "Chile Primera Division" sheet were re-ordered:
# the id-column (A) stays put, but the rest of each affected row's data
# (columns B..AD -- match id, teams, scores, odds, etc.) is rotated among
# a handful of rows so each match's full record lands on the correct row.
#
# Cycles (old row -> new row that receives its data):
#   105 -> 106, 106 -> 105
#   118 -> 119, 119 -> 120, 120 -> 118
#   121 -> 122, 122 -> 121
#   137 -> 138, 138 -> 137
#   139 -> 140, 140 -> 139

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

# Each inner array is one rotation cycle, listed in the order that data
# flows: row[0]'s original data ends up on row[1], row[1]'s on row[2], ...
# and the last row's data wraps back around to row[0].
$cycles = @(
    @(106, 105),
    @(120, 119, 118),
    @(122, 121),
    @(138, 137),
    @(140, 139)
)

foreach ($cycle in $cycles) {
    # Snapshot B..AD for every row in this cycle before writing anything.
    $snapshots = @{}
    foreach ($row in $cycle) {
        $rowValues = @{}
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $rowValues[$col] = $ws.Cells.Item($row, $col).Value()
        }
        $snapshots[$row] = $rowValues
    }

    # Write row[i]'s snapshot into row[i+1] (wrapping around).
    $count = $cycle.Count
    for ($i = 0; $i -lt $count; $i++) {
        $srcRow = $cycle[$i]
        $dstRow = $cycle[($i + 1) % $count]
        $rowValues = $snapshots[$srcRow]
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $ws.Cells.Item($dstRow, $col).Value = $rowValues[$col]
        }
    }
}
